{"js": "// Apply the resume redaction/update edit described by the diff:\n// 1. Name -> title case\n// 2. Professional title -> generic placeholder\n// 3. Phone/email formatting -> generic style\n// 4. \"20+ years\" -> \"21 years\"\n// 5. Company name -> generic placeholder\n// 6. Remove 5 whole job entries (DATA PRODUCTS MANAGER .. RESEARCH DIRECTOR)\n// 7. Remove product/codename call-outs (BALLISTA, DAMON, SimCrisis, RACSO) from\n//    the \"Software Development and Innovation\" bullets\n// 8. Remove the \"Data Architecture and Engineering\" and\n//    \"Research Impact and Recognition\" achievement sections entirely\n\nconst body = context.document.body;\n\n// --- 1) Simple literal text replacements -------------------------------\nasync function replaceOnce(searchText, replacement, matchCase) {\n  const results = body.search(searchText, { matchCase: matchCase !== false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nawait replaceOnce(\"DHEERAJ CHAND\", \"Dheeraj Chand\");\nawait replaceOnce(\n  \"Research, Data Analytics & Engineering Professional\",\n  \"Professional Title\"\n);\nawait replaceOnce(\n  \"(202) 550-7110 | Dheeraj.Chand@gmail.com\",\n  \"202.550.7110 | dheeraj.chand@gmail.com\"\n);\nawait replaceOnce(\n  \"Research & Data Professional with 20+ years of comprehensive experience spanning applied research, data engineering, and software development. Expert in translating complex analytical requirements into scalable technical solutions. Proven track record leading cross-functional teams, architecting data platforms, and delivering insights that drive strategic decision-making across political, nonprofit, and technology sectors. Deep expertise in survey methodology, geospatial analysis, and building production systems for sensitive data applications.\",\n  \"Research & Data Professional with 21 years of comprehensive experience spanning applied research, data engineering, and software development. Expert in translating complex analytical requirements into scalable technical solutions. Proven track record leading cross-functional teams, architecting data platforms, and delivering insights that drive strategic decision-making across political, nonprofit, and technology sectors. Deep expertise in survey methodology, geospatial analysis, and building production systems for sensitive data applications.\"\n);\nawait replaceOnce(\n  \"Siege Analytics, Austin, TX | 2005 \u2013 Present\",\n  \"Your Company Name, Your City, ST | 2005 \u2013 Present\"\n);\n\nawait replaceOnce(\n  \"\u2713 Conceived and deployed BALLISTA redistricting software used by thousands of analysts nationwide\",\n  \"\u2713 Conceived and deployed redistricting software used by thousands of analysts nationwide\"\n);\nawait replaceOnce(\n  \"\u2713 Developed DAMON boundary estimation system using incomplete data without ML requirements\",\n  \"\u2713 Developed boundary estimation system using incomplete data without ML requirements\"\n);\nawait replaceOnce(\n  \"\u2713 Created SimCrisis econometric simulation platform for humanitarian intervention modeling\",\n  \"\u2713 Created econometric simulation platform for humanitarian intervention modeling\"\n);\nawait replaceOnce(\n  \"\u2713 Built RACSO comprehensive survey operations platform from RFP through deployment\",\n  \"\u2713 Built comprehensive survey operations platform from RFP through deployment\"\n);\n\n// --- 2) Remove whole paragraph ranges -----------------------------------\n// Re-load the (now updated) paragraph list and delete by text-anchored\n// boundaries so the logic survives any paragraph-count drift caused by the\n// edits above.\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nfunction indexOfText(items, text) {\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text === text) {\n      return i;\n    }\n  }\n  return -1;\n}\n\nconst items = paras.items;\n\n// Block 1: from \"DATA PRODUCTS MANAGER\" heading through the last bullet of\n// the RESEARCH DIRECTOR entry (\"Created comprehensive data visualization\n// solutions for complex research findings\"), inclusive.\nconst block1Start = indexOfText(items, \"DATA PRODUCTS MANAGER\");\nconst block1End = indexOfText(\n  items,\n  \"\u25b8 Created comprehensive data visualization solutions for complex research findings\"\n);\n\n// Block 2: from \"Data Architecture and Engineering\" heading through the\n// last bullet of \"Research Impact and Recognition\"\n// (\"\u2713 Pioneered integration of geospatial techniques into political and\n// market research\"), inclusive.\nconst block2Start = indexOfText(items, \"Data Architecture and Engineering\");\nconst block2End = indexOfText(\n  items,\n  \"\u2713 Pioneered integration of geospatial techniques into political and market research\"\n);\n\nconst ranges = [];\nif (block1Start !== -1 && block1End !== -1 && block1End >= block1Start) {\n  ranges.push([block1Start, block1End]);\n}\nif (block2Start !== -1 && block2End !== -1 && block2End >= block2Start) {\n  ranges.push([block2Start, block2End]);\n}\n\n// Delete from the bottom up so earlier indices stay valid.\nranges.sort((a, b) => b[0] - a[0]);\nfor (const [start, end] of ranges) {\n  for (let i = end; i >= start; i--) {\n    items[i].delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply the resume redaction/update edit described by the diff:\n# 1. Name -> title case\n# 2. Professional title -> generic placeholder\n# 3. Phone/email formatting -> generic style\n# 4. \"20+ years\" -> \"21 years\"\n# 5. Company name -> generic placeholder\n# 6. Remove 5 whole job entries (DATA PRODUCTS MANAGER .. RESEARCH DIRECTOR)\n# 7. Remove product/codename call-outs (BALLISTA, DAMON, SimCrisis, RACSO) from\n#    the \"Software Development and Innovation\" bullets\n# 8. Remove the \"Data Architecture and Engineering\" and\n#    \"Research Impact and Recognition\" achievement sections entirely\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText, $matchCase) {\n    $find = $d.Content.Find\n    $find.Execute($findText, $matchCase, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\nReplace-Text \"DHEERAJ CHAND\" \"Dheeraj Chand\" $true\nReplace-Text \"Research, Data Analytics & Engineering Professional\" \"Professional Title\" $true\nReplace-Text \"(202) 550-7110 | Dheeraj.Chand@gmail.com\" \"202.550.7110 | dheeraj.chand@gmail.com\" $true\nReplace-Text \"Research & Data Professional with 20+ years of comprehensive experience spanning applied research, data engineering, and software development. Expert in translating complex analytical requirements into scalable technical solutions. Proven track record leading cross-functional teams, architecting data platforms, and delivering insights that drive strategic decision-making across political, nonprofit, and technology sectors. Deep expertise in survey methodology, geospatial analysis, and building production systems for sensitive data applications.\" \"Research & Data Professional with 21 years of comprehensive experience spanning applied research, data engineering, and software development. Expert in translating complex analytical requirements into scalable technical solutions. Proven track record leading cross-functional teams, architecting data platforms, and delivering insights that drive strategic decision-making across political, nonprofit, and technology sectors. Deep expertise in survey methodology, geospatial analysis, and building production systems for sensitive data applications.\" $true\nReplace-Text \"Siege Analytics, Austin, TX | 2005 \u2013 Present\" \"Your Company Name, Your City, ST | 2005 \u2013 Present\" $true\n\nReplace-Text \"\u2713 Conceived and deployed BALLISTA redistricting software used by thousands of analysts nationwide\" \"\u2713 Conceived and deployed redistricting software used by thousands of analysts nationwide\" $true\nReplace-Text \"\u2713 Developed DAMON boundary estimation system using incomplete data without ML requirements\" \"\u2713 Developed boundary estimation system using incomplete data without ML requirements\" $true\nReplace-Text \"\u2713 Created SimCrisis econometric simulation platform for humanitarian intervention modeling\" \"\u2713 Created econometric simulation platform for humanitarian intervention modeling\" $true\nReplace-Text \"\u2713 Built RACSO comprehensive survey operations platform from RFP through deployment\" \"\u2713 Built comprehensive survey operations platform from RFP through deployment\" $true\n\n# --- Remove whole paragraph ranges, matched by exact paragraph text so the\n# logic is resilient to any paragraph-count drift caused by the edits above.\nfunction Find-ParaIndex($target) {\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        $t = $d.Paragraphs.Item($i).Range.Text\n        $t = $t.TrimEnd([char]13)\n        if ($t -eq $target) {\n            return $i\n        }\n    }\n    return -1\n}\n\nfunction Delete-ParaRange($startText, $endText) {\n    $startIdx = Find-ParaIndex $startText\n    $endIdx = Find-ParaIndex $endText\n    if ($startIdx -ne -1 -and $endIdx -ne -1 -and $endIdx -ge $startIdx) {\n        $start = $d.Paragraphs.Item($startIdx).Range.Start\n        $end = $d.Paragraphs.Item($endIdx).Range.End\n        $r = $d.Range($start, $end)\n        $r.Delete()\n    }\n}\n\n# Block 1: DATA PRODUCTS MANAGER .. last bullet of RESEARCH DIRECTOR\nDelete-ParaRange \"DATA PRODUCTS MANAGER\" \"\u25b8 Created comprehensive data visualization solutions for complex research findings\"\n\n# Block 2: Data Architecture and Engineering .. last bullet of Research Impact and Recognition\nDelete-ParaRange \"Data Architecture and Engineering\" \"\u2713 Pioneered integration of geospatial techniques into political and market research\"\n"}
